$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Write-Host ("Sheet name: " + $ws.Name)
$v1 = $ws.Range("A7").Value()
Write-Host ("A7: " + $v1)
$v2 = $ws.Range("A9").Value2()
Write-Host ("A9: " + $v2)
$v3 = $ws.Cells.Item(7,1).Value()
Write-Host ("Cells(7,1): " + $v3)
